# Replace three-digit-by-one-digit multiplication problems with new values
$d = $word.ActiveDocument

$replacements = @(
    @("658×6=3948", "852×2=1704"),
    @("799×6=4794", "544×3=1632"),
    @("850×2=1700", "225×2=450"),
    @("398×9=3582", "622×8=4976"),
    @("367×6=2202", "838×8=6704"),
    @("867×4=3468", "219×9=1971"),
    @("261×2=522",  "448×6=2688"),
    @("749×4=2996", "364×7=2548"),
    @("380×2=760",  "399×7=2793"),
    @("326×4=1304", "904×2=1808"),
    @("817×7=5719", "669×6=4014"),
    @("301×3=903",  "128×4=512"),
    @("528×5=2640", "169×7=1183"),
    @("460×8=3680", "143×4=572"),
    @("837×3=2511", "620×7=4340"),
    @("429×2=858",  "619×6=3714"),
    @("893×3=2679", "305×7=2135"),
    @("741×7=5187", "321×9=2889"),
    @("269×3=807",  "857×9=7713"),
    @("744×2=1488", "516×2=1032"),
    @("335×3=1005", "522×7=3654"),
    @("525×2=1050", "265×8=2120"),
    @("490×4=1960", "182×6=1092"),
    @("269×7=1883", "182×2=364"),
    @("223×8=1784", "878×9=7902")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
